# Update country data files
# Insert the "Number of employees / Assets / Turnover" MSME size-definition
# table above the existing "Sector Distribution Details" table, and append
# the INEGI source-attribution rows at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 6 blank rows at row 24. Everything from the old
#    row 24 ("Sector Distribution Details") down to the old row 38
#    shifts down by 6 rows (24->30 ... 38->44), which is exactly what
#    the target layout needs. Existing cell values/styles move with
#    their rows automatically.
# ---------------------------------------------------------------------
$ws.Rows("24:29").Insert()

# ---------------------------------------------------------------------
# 2. Fill in the new "Number of employees / Assets / Turnover" table in
#    the newly freed rows 21-25.
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "Number of employees"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Font.Bold = $true

$ws.Range("A22").Value = "Micro"
$ws.Range("B22").Value = "0-10"

$ws.Range("A23").Value = "Small"
$ws.Range("B23").Value = "10-50"

$ws.Range("A24").Value = "Medium"
$ws.Range("B24").Value = "50-250"

$ws.Range("A25").Value = "Large"
$ws.Range("B25").Value = ">250"

# ---------------------------------------------------------------------
# 3. Fix up the source-attribution block (now rows 37-40) and append the
#    new INEGI citation rows (43-44).
# ---------------------------------------------------------------------
$ws.Range("A43").Value = "INEGI"
$ws.Range("A43").Font.Bold = $true

$ws.Range("A44").Value = "Instituto Nacional de Estadistica y Geografia (INEGI), ""Resumen de los resultados de los Censos Económicos 2009; Micro, pequeña, mediana, y gran empresa. Estratificación de los establecimientos"", 2011, p.17-18. Available at http://www.inegi.org.mx/prod_serv/contenidos/espanol/bvinegi/productos/censos/economicos/2009/comercio/micro_peque_media/Mono_Micro_peque_mediana.pdf"
$ws.Range("A44").Font.Italic = $true

# ---------------------------------------------------------------------
# 4. The row-insert does not move the worksheet's <hyperlinks> entry, so
#    the hyperlink on the INEGI PDF URL cell (old A33, now A39) has to be
#    re-pointed by hand.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A39"), "http://www.inegi.org.mx/prod_serv/contenidos/espanol/bvinegi/productos/censos/economicos/2009/comercio/micro_peque_media/Mono_Micro_peque_mediana.pdf")
